$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Shared text updates (Volume number, week-of dates) ---
$ws.Range("A8").Value2 = "Volume 30   Number  3"
$ws.Range("C9").Value2 = "Report Covering the Week  1/16/2023  Through  1/22/2023"

# --- Cells that change type/style (copy formatting from a stable template cell, then set the value) ---
$ws.Range("C30").Copy($ws.Range("D14"))
$ws.Range("E30").Copy($ws.Range("E14"))
$ws.Range("K36").Copy($ws.Range("L15"))
$ws.Range("L15").Value2 = 100
$ws.Range("C30").Copy($ws.Range("D22"))
$ws.Range("E30").Copy($ws.Range("E22"))
$ws.Range("C30").Copy($ws.Range("C23"))
$ws.Range("C36").Copy($ws.Range("D23"))
$ws.Range("D23").Value2 = 1
$ws.Range("K36").Copy($ws.Range("E23"))
$ws.Range("E23").Value2 = -100
$ws.Range("C36").Copy($ws.Range("G23"))
$ws.Range("G23").Value2 = 1
$ws.Range("K36").Copy($ws.Range("H23"))
$ws.Range("H23").Value2 = 200
$ws.Range("C36").Copy($ws.Range("J23"))
$ws.Range("J23").Value2 = 1
$ws.Range("K36").Copy($ws.Range("K23"))
$ws.Range("K23").Value2 = 100
$ws.Range("K36").Copy($ws.Range("L23"))
$ws.Range("L23").Value2 = 100
$ws.Range("K36").Copy($ws.Range("M23"))
$ws.Range("M23").Value2 = 100
$ws.Range("C36").Copy($ws.Range("C26"))
$ws.Range("C26").Value2 = 1
$ws.Range("K36").Copy($ws.Range("L26"))
$ws.Range("L26").Value2 = 100
$ws.Range("C36").Copy($ws.Range("D27"))
$ws.Range("D27").Value2 = 2
$ws.Range("K36").Copy($ws.Range("E27"))
$ws.Range("E27").Value2 = -50
$ws.Range("C30").Copy($ws.Range("C28"))
$ws.Range("C30").Copy($ws.Range("D28"))
$ws.Range("E30").Copy($ws.Range("E28"))
$ws.Range("C30").Copy($ws.Range("C29"))
$ws.Range("C30").Copy($ws.Range("D29"))
$ws.Range("E30").Copy($ws.Range("E29"))

# --- Plain numeric value updates ---
$ws.Range("F15").Value2 = 2
$ws.Range("G15").Value2 = 2
$ws.Range("H15").Value2 = 0
$ws.Range("J15").Value2 = 2
$ws.Range("K15").Value2 = 0
$ws.Range("N15").Value2 = -50
$ws.Range("F16").Value2 = 23
$ws.Range("G16").Value2 = 13
$ws.Range("H16").Value2 = 76.923076923076
$ws.Range("I16").Value2 = 19
$ws.Range("J16").Value2 = 11
$ws.Range("K16").Value2 = 72.727272727272
$ws.Range("L16").Value2 = 26.666666666666
$ws.Range("M16").Value2 = 72.727272727272
$ws.Range("N16").Value2 = -76.25
$ws.Range("C17").Value2 = 12
$ws.Range("D17").Value2 = 3
$ws.Range("E17").Value2 = 300
$ws.Range("F17").Value2 = 36
$ws.Range("G17").Value2 = 12
$ws.Range("H17").Value2 = 200
$ws.Range("I17").Value2 = 28
$ws.Range("J17").Value2 = 9
$ws.Range("K17").Value2 = 211.111111111111
$ws.Range("L17").Value2 = 133.333333333333
$ws.Range("M17").Value2 = 211.111111111111
$ws.Range("N17").Value2 = 47.368421052631
$ws.Range("C18").Value2 = 16
$ws.Range("D18").Value2 = 12
$ws.Range("E18").Value2 = 33.333333333333
$ws.Range("F18").Value2 = 54
$ws.Range("G18").Value2 = 38
$ws.Range("H18").Value2 = 42.105263157894
$ws.Range("I18").Value2 = 41
$ws.Range("J18").Value2 = 28
$ws.Range("K18").Value2 = 46.428571428571
$ws.Range("L18").Value2 = 51.851851851851
$ws.Range("M18").Value2 = 32.258064516129
$ws.Range("N18").Value2 = -76.966292134831
$ws.Range("C19").Value2 = 27
$ws.Range("D19").Value2 = 39
$ws.Range("E19").Value2 = -30.769230769230
$ws.Range("F19").Value2 = 110
$ws.Range("G19").Value2 = 165
$ws.Range("H19").Value2 = -33.333333333333
$ws.Range("I19").Value2 = 92
$ws.Range("J19").Value2 = 116
$ws.Range("K19").Value2 = -20.689655172413
$ws.Range("L19").Value2 = 206.666666666667
$ws.Range("M19").Value2 = 113.953488372093
$ws.Range("N19").Value2 = 21.052631578947
$ws.Range("C20").Value2 = 4
$ws.Range("D20").Value2 = 8
$ws.Range("E20").Value2 = -50
$ws.Range("F20").Value2 = 37
$ws.Range("G20").Value2 = 22
$ws.Range("H20").Value2 = 68.181818181818
$ws.Range("I20").Value2 = 23
$ws.Range("J20").Value2 = 19
$ws.Range("K20").Value2 = 21.052631578947
$ws.Range("L20").Value2 = 76.923076923076
$ws.Range("M20").Value2 = 27.777777777777
$ws.Range("N20").Value2 = -91.385767790262
$ws.Range("C21").Value2 = 67
$ws.Range("D21").Value2 = 67
$ws.Range("E21").Value2 = 0
$ws.Range("F21").Value2 = 262
$ws.Range("G21").Value2 = 253
$ws.Range("H21").Value2 = 3.557312252964
$ws.Range("I21").Value2 = 205
$ws.Range("J21").Value2 = 186
$ws.Range("K21").Value2 = 10.215053763440
$ws.Range("L21").Value2 = 109.183673469388
$ws.Range("M21").Value2 = 79.824561403508
$ws.Range("N21").Value2 = -67.2
$ws.Range("F22").Value2 = 3
$ws.Range("G22").Value2 = 2
$ws.Range("H22").Value2 = 50
$ws.Range("C24").Value2 = 57
$ws.Range("D24").Value2 = 56
$ws.Range("E24").Value2 = 1.785714285714
$ws.Range("F24").Value2 = 204
$ws.Range("H24").Value2 = -2.857142857142
$ws.Range("I24").Value2 = 166
$ws.Range("J24").Value2 = 161
$ws.Range("K24").Value2 = 3.105590062111
$ws.Range("L24").Value2 = 48.214285714285
$ws.Range("M24").Value2 = 62.745098039215
$ws.Range("C25").Value2 = 15
$ws.Range("D25").Value2 = 9
$ws.Range("E25").Value2 = 66.666666666666
$ws.Range("F25").Value2 = 55
$ws.Range("G25").Value2 = 45
$ws.Range("H25").Value2 = 22.222222222222
$ws.Range("I25").Value2 = 44
$ws.Range("J25").Value2 = 34
$ws.Range("K25").Value2 = 29.411764705882
$ws.Range("L25").Value2 = 109.52380952381
$ws.Range("M25").Value2 = 41.935483870967
$ws.Range("D26").Value2 = 1
$ws.Range("E26").Value2 = 0
$ws.Range("F26").Value2 = 5
$ws.Range("G26").Value2 = 3
$ws.Range("H26").Value2 = 66.666666666666
$ws.Range("I26").Value2 = 4
$ws.Range("J26").Value2 = 3
$ws.Range("K26").Value2 = 33.333333333333
$ws.Range("C27").Value2 = 1
$ws.Range("F27").Value2 = 3
$ws.Range("G27").Value2 = 6
$ws.Range("I27").Value2 = 3
$ws.Range("J27").Value2 = 5
$ws.Range("K27").Value2 = -40
$ws.Range("L27").Value2 = 0
